# C5-PowerPoint.pptx edit
#
# 1) Slide 6, shape 2 is a table (graphicFrame) whose table style is
#    switched from the default "Table_0" style ({DEF32ECF-0DE1-4CDC-A063-
#    DB070FDE6E2C}) to a built-in themed table style
#    ({616AACAD-8BAF-4749-97A4-DB3BD5873582}).
#
# 2) The deck's design theme is changed from the "Integral" Office theme
#    to the plain default "Office Theme" - i.e. the slide master's
#    12-colour theme palette is swapped from the Integral colours to the
#    stock Office colours.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 6 -------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{616AACAD-8BAF-4749-97A4-DB3BD5873582}")

# --- 2. Switch the slide master's theme colours to the default Office ----
#        palette (was the "Integral" theme's palette).
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

$colorScheme.Item(1).RGB  = 0         # dk1      000000
$colorScheme.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388   # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501   # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407     # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308  # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456   # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797  # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477   # folHlink 954F72
